$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column F, matching the style of the existing header row (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Populate the time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 13:40:38.220211"
$ws.Range("F3").Value = "2021-10-05 13:40:38.220224"
$ws.Range("F4").Value = "2021-10-05 13:40:38.220228"
$ws.Range("F5").Value = "2021-10-05 13:40:38.220231"
$ws.Range("F6").Value = "2021-10-05 13:40:38.220234"
$ws.Range("F7").Value = "2021-10-05 13:40:38.220237"
$ws.Range("F8").Value = "2021-10-05 13:40:38.220240"
